# Insert a new weekly price record for "Macroferia Regional de Talca - Cilantro".
# This pushes the existing data rows 39-86 down to 40-87 and places the new
# observation (fecha 45049 = 2023-05-03) into the now-empty row 39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 39, shifting rows 39:86 down to 40:87.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 45049
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 100112040
$ws.Range("G39").Value = "Cilantro"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 150
$ws.Range("K39").Value = 7000
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = 7000
$ws.Range("N39").Value = "$/caja 36 atados"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 194
$ws.Range("Q39").Value = 36
$ws.Range("R39").Value = "Hortaliza"
